# Refresh the cryptocurrency price/volume snapshot (columns D and E)
# for the rows whose figures moved since the last GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.470.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "'1.904.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'325.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.61%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "'0.4792"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.42%  "
$ws.Range("D8").Value = "'0.4065"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("D9").Value = "'0.08073"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("D11").Value = "'23.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.39%  "
$ws.Range("D12").Value = "'1.899.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").Value = "'5.951"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("D15").Value = "'89.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").Value = "'0.06707"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.84%  "
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").Value = "'17.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "'29.474.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").Value = "'5.539"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").Value = "'2.165"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.14%  "
$ws.Range("D25").Value = "'2.124.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").Value = "'154.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("D27").Value = "'19.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").Value = "'6.095"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.98%  "
$ws.Range("D29").Value = "'2.096"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("D30").Value = "'118.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("D31").Value = "'1.038"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.96%  "
$ws.Range("D32").Value = "'0.09509"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D33").Value = "'5.436"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("E34").Value = "  -2.55%  "
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").Value = "'0.06076"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("D38").Value = "'1.176"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("D39").Value = "'0.5877"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("D40").Value = "'7.919"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.92%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "'10.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").Value = "'2.422"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.30%  "
$ws.Range("D44").Value = "'1.279"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("D45").Value = "'0.07791"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.61%  "
$ws.Range("D46").Value = "'12.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("D47").Value = "'0.5528"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("D48").Value = "'1.921"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("D49").Value = "'113.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").Value = "'0.2941"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("E51").Value = "  +0.92%  "
